$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Nourhan Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy, Dr. Rana Abo-Zaid, Dr. Shimaa Ahmad Mekki"
$ws.Range("G4").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G5").Value = "Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G6").Value = "Dr. Kerelos Zareef, Dr. Nada Mohammad"
$ws.Range("G10").Value = "Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Madeha Saeed"
$ws.Range("G15").Value = "Dr. Hana Amr, Dr. Nancy Abd Al-Shafy"
$ws.Range("G18").Value = "Dr. Yasmin, Dr. Aya Hanafy, Dr. Remon, Dr. Shorok Mohammad"
$ws.Range("G19").Value = "Dr. Remon, Dr. Nardine, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Monica, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Yassmen Ahmad"
$ws.Range("G20").Value = "Dr. Remon, Dr. Nardine, Dr. Youstina Magdy, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Marina Sorial"
$ws.Range("G21").Value = "Dr. Yasmin, Dr. Monica, Dr. Shorok Mohammad, Dr. Neveen Nashaat, Dr. Yassmen Ahmad"
$ws.Range("G22").Value = "Dr. Remon, Dr. Monica, Dr. Wafaa Ebida, Dr. Naema Gomaa"
$ws.Range("G23").Value = "Dr. Yassmen Ahmad, Dr. Wafaa Ebida"
$ws.Range("G24").Value = "Dr. Nourhan Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy, Dr. Rana Abo-Zaid, Dr. Shimaa Ahmad Mekki"
$ws.Range("G25").Value = "Administrator, Dr. Gehan Adel, Dr. Alshimaa Atef, Dr. Manar Montaser"
$ws.Range("G26").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy"
$ws.Range("G27").Value = "Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G28").Value = "Dr. Kerelos Zareef, Dr. Nada Mohammad"
$ws.Range("G32").Value = "Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Madeha Saeed"
$ws.Range("G37").Value = "Dr. Hana Amr, Dr. Nancy Abd Al-Shafy"
$ws.Range("G40").Value = "Dr. Yasmin, Dr. Aya Hanafy, Dr. Remon, Dr. Shorok Mohammad"
$ws.Range("G41").Value = "Dr. Remon, Dr. Nardine, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Monica, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Yassmen Ahmad"
$ws.Range("G42").Value = "Dr. Remon, Dr. Nardine, Dr. Youstina Magdy, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Marina Sorial"
$ws.Range("G43").Value = "Dr. Yasmin, Dr. Monica, Dr. Shorok Mohammad, Dr. Neveen Nashaat, Dr. Yassmen Ahmad"
$ws.Range("G44").Value = "Dr. Remon, Dr. Monica, Dr. Wafaa Ebida, Dr. Naema Gomaa"
$ws.Range("G45").Value = "Dr. Yassmen Ahmad, Dr. Wafaa Ebida"
$ws.Range("G46").Value = "Dr. Nahla Nagiub, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud"
$ws.Range("G48").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Nahla Nagiub, Dr. Servinaz Sayed Mohammad"
$ws.Range("G49").Value = "Dr. Amera Ahmad Saad, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad"
$ws.Range("G52").Value = "Dr. Shimaa Ashraf, Dr. Mariam Nour El-Din"
$ws.Range("G54").Value = "Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Merna Said, Dr. Mai Mustafa, Dr. Arwa Al-Sayed, Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G59").Value = "Dr. Walaa Ghanima, Dr. Marian Samir, Dr. Enas Omran"
$ws.Range("G60").Value = "Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy"
$ws.Range("G62").Value = "Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Shorok Mohammad"
$ws.Range("G63").Value = "Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Neveen Nashaat"
$ws.Range("G64").Value = "Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Youstina Magdy"
$ws.Range("G65").Value = "Dr. Eman Samir Gabry, Dr. Remon, Dr. Nardine, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Shorok Mohammad, Dr. Neveen Nashaat"
$ws.Range("G66").Value = "Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Maryam Ashraf, Dr. Monica, Dr. Marina Sorial"
$ws.Range("G68").Value = "Dr. Nahla Nagiub, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Hend Mahmoud"
$ws.Range("G70").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Nahla Nagiub, Dr. Servinaz Sayed Mohammad"
$ws.Range("G71").Value = "Dr. Amera Ahmad Saad, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad"
$ws.Range("G74").Value = "Dr. Shimaa Ashraf, Dr. Mariam Nour El-Din"
$ws.Range("G76").Value = "Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Merna Said, Dr. Mai Mustafa, Dr. Arwa Al-Sayed, Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G81").Value = "Dr. Walaa Ghanima, Dr. Marian Samir, Dr. Enas Omran"
$ws.Range("G82").Value = "Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy"
$ws.Range("G84").Value = "Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Shorok Mohammad"
$ws.Range("G85").Value = "Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Neveen Nashaat"
$ws.Range("G86").Value = "Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Youstina Magdy"
$ws.Range("G87").Value = "Dr. Eman Samir Gabry, Dr. Remon, Dr. Nardine, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Shorok Mohammad, Dr. Neveen Nashaat"
$ws.Range("G88").Value = "Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Maryam Ashraf, Dr. Monica, Dr. Marina Sorial"
$ws.Range("G90").Value = "Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser"
$ws.Range("G92").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Nahla Nagiub, Dr. Servinaz Sayed Mohammad"
$ws.Range("G93").Value = "Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G96").Value = "Dr. Nourhan Mohammad, Dr. Sara Nabil, Dr. Mariam Nour El-Din, Dr. Amal Awwad"
$ws.Range("G98").Value = "Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Merna Said, Dr. Mai Mustafa, Dr. Arwa Al-Sayed, Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G104").Value = "Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy"
$ws.Range("G106").Value = "Dr. Remon, Dr. Youstina Magdy, Dr. Nardine, Dr. Wafaa Ebida, Dr. Monica, Dr. Neveen Nashaat"
$ws.Range("G107").Value = "Dr. Aya Hanafy, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Monica, Dr. Neveen Nashaat, Dr. Yassmen Ahmad"
$ws.Range("G108").Value = "Dr. Remon, Dr. Nardine, Dr. Youstina Magdy, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Marina Sorial"
$ws.Range("G110").Value = "Dr. Monica, Dr. Yassmen Ahmad, Dr. Wafaa Ebida"
$ws.Range("G111").Value = "Dr. Eman Samir Gabry, Dr. Marina Atef, Dr. Yasmin, Dr. Nourham Mostafa, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Monica"
$ws.Range("G112").Value = "Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser"
$ws.Range("G114").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Nahla Nagiub, Dr. Servinaz Sayed Mohammad"
$ws.Range("G115").Value = "Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G118").Value = "Dr. Nourhan Mohammad, Dr. Sara Nabil, Dr. Mariam Nour El-Din, Dr. Amal Awwad"
$ws.Range("G120").Value = "Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Merna Said, Dr. Mai Mustafa, Dr. Arwa Al-Sayed, Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G126").Value = "Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy"
$ws.Range("G128").Value = "Dr. Remon, Dr. Youstina Magdy, Dr. Nardine, Dr. Wafaa Ebida, Dr. Monica, Dr. Neveen Nashaat"
$ws.Range("G129").Value = "Dr. Aya Hanafy, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Monica, Dr. Neveen Nashaat, Dr. Yassmen Ahmad"
$ws.Range("G130").Value = "Dr. Remon, Dr. Nardine, Dr. Youstina Magdy, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Marina Sorial"
$ws.Range("G132").Value = "Dr. Monica, Dr. Yassmen Ahmad, Dr. Wafaa Ebida"
$ws.Range("G133").Value = "Dr. Eman Samir Gabry, Dr. Marina Atef, Dr. Yasmin, Dr. Nourham Mostafa, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Monica"
$ws.Range("G134").Value = "Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G137").Value = "Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G142").Value = "Dr. Basma Hamed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh"
$ws.Range("G147").Value = "Dr. Nourham Mostafa, Dr. Nancy Abd Al-Shafy"
$ws.Range("G150").Value = "Dr. Remon, Dr. Nardine, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Monica, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Yassmen Ahmad"
$ws.Range("G151").Value = "Dr. Yassmen Ahmad, Dr. Monica, Dr. Marina Atef, Dr. Wafaa Ebida"
$ws.Range("G152").Value = "Dr. Marina Atef, Dr. Wafaa Ebida"
$ws.Range("G153").Value = "Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Maryam Ashraf, Dr. Monica, Dr. Marina Sorial"
$ws.Range("G154").Value = "Dr. Remon, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Naema Gomaa"
$ws.Range("G155").Value = "Dr. Eman Samir Gabry, Dr. Marina Atef, Dr. Yasmin, Dr. Nourham Mostafa, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Monica"
$ws.Range("G156").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany"
$ws.Range("G159").Value = "Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G164").Value = "Dr. Basma Hamed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh"
$ws.Range("G169").Value = "Dr. Nourham Mostafa, Dr. Nancy Abd Al-Shafy"
$ws.Range("G172").Value = "Dr. Remon, Dr. Nardine, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Monica, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Yassmen Ahmad"
$ws.Range("G173").Value = "Dr. Yassmen Ahmad, Dr. Monica, Dr. Marina Atef, Dr. Wafaa Ebida"
$ws.Range("G174").Value = "Dr. Marina Atef, Dr. Wafaa Ebida"
$ws.Range("G175").Value = "Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Maryam Ashraf, Dr. Monica, Dr. Marina Sorial"
$ws.Range("G176").Value = "Dr. Remon, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Naema Gomaa"
$ws.Range("G177").Value = "Dr. Eman Samir Gabry, Dr. Marina Atef, Dr. Yasmin, Dr. Nourham Mostafa, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Monica"
